# Update "想去人数" (interested count) figures for the same conan/expo
# entries on both the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 715
    $ws.Range("F3").Value = 4060
}
